# Updates cryptos list values (price column D, 1h volume % column E)
# to match the refreshed scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Cells whose new text would otherwise be auto-parsed as a number by
# Excel get an explicit Text number format first so the literal digits
# (including trailing zeros / decimal grouping) are preserved verbatim.
$ws.Range("D2").Value = "37.818.25"
$ws.Range("D3").Value = "2.088.84"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.03"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.65"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.18"
$ws.Range("D13").Value = "2.396.43"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.36"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.781"
$ws.Range("D17").Value = "2.093.44"
$ws.Range("D18").Value = "37.782.25"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.29"
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.39"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.81"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.84"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.50"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.57"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0975"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.94"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.90"
$ws.Range("D45").Value = "1.453.55"
$ws.Range("D51").Value = "2.280.63"

# --- Column E (Volume 1h %) updates ---
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  +8.96%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  +9.66%  "
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  +0.04%  "

